# Update gh-pages output data (想去人数 / "interested" counts) in the
# 展览 (Exhibitions) and 全部类型 (All Types) sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 5170
$wsExhibit.Range("F7").Value = 60
$wsExhibit.Range("F9").Value = 339

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value = 5170
$wsAll.Range("F11").Value = 60
$wsAll.Range("F14").Value = 339
